$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 69.28570999999999
$ws.Range("I5").Value = 69.28570999999999
$ws.Range("K5").Value = 69.28570999999999
$ws.Range("M5").Value = 45.71429000000001
$ws.Range("H43").Value = 261688.62
$ws.Range("I43").Value = 9500
$ws.Range("J43").Value = 345751.5
$ws.Range("K43").Value = 9500
$ws.Range("L43").Value = 345751.5
$ws.Range("M43").Value = -9431
$ws.Range("N43").Value = -345889.5
$ws.Range("H129").Value = 1960.2222
$ws.Range("I129").Value = 955.25
$ws.Range("K129").Value = 2865.75
$ws.Range("M129").Value = 2134.25
$ws.Range("H135").Value = 1404.1364
$ws.Range("I135").Value = 1535.1765
$ws.Range("K135").Value = 13816.5885
$ws.Range("M135").Value = -11281.5885
$ws.Range("H137").Value = 4533.579
$ws.Range("I137").Value = 2458.5
$ws.Range("K137").Value = 7375.5
$ws.Range("M137").Value = -4825.5
$ws.Range("H138").Value = 2809.4043
$ws.Range("J138").Value = 3270.676
$ws.Range("L138").Value = 9812.028
$ws.Range("N138").Value = -20092.028
$ws.Range("H141").Value = 1557.4231
$ws.Range("I141").Value = 1557.4231
$ws.Range("K141").Value = 4672.2693
$ws.Range("M141").Value = 507.7307000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 20053624
$ws.Range("I61").Value = 25012944
$ws.Range("K61").Value = 25012944
$ws.Range("M61").Value = -25012732
$ws.Range("H69").Value = 194996.67
$ws.Range("J69").Value = 194996.67
$ws.Range("L69").Value = 194996.67
$ws.Range("N69").Value = -196494.67
$ws.Range("H72").Value = 194996.67
$ws.Range("J72").Value = 194996.67
$ws.Range("L72").Value = 584990.01
$ws.Range("N72").Value = -592478.01
$ws.Range("H88").Value = 2907.9092
$ws.Range("I88").Value = 3780
$ws.Range("J88").Value = 2181.1667
$ws.Range("K88").Value = 3780
$ws.Range("L88").Value = 2181.1667
$ws.Range("M88").Value = -3374
$ws.Range("N88").Value = -2993.1667
$ws.Range("H91").Value = 2907.9092
$ws.Range("I91").Value = 3780
$ws.Range("J91").Value = 2181.1667
$ws.Range("K91").Value = 3780
$ws.Range("L91").Value = 2181.1667
$ws.Range("M91").Value = -2376
$ws.Range("N91").Value = -4989.1667
$ws.Range("H97").Value = 1637.6666
$ws.Range("I97").Value = 1720.6154
$ws.Range("K97").Value = 1720.6154
$ws.Range("M97").Value = -1224.6154
$ws.Range("H122").Value = 2017.125
$ws.Range("I122").Value = 2133.8572
$ws.Range("J122").Value = 1200
$ws.Range("K122").Value = 6401.571599999999
$ws.Range("L122").Value = 3600
$ws.Range("M122").Value = -3951.571599999999
$ws.Range("N122").Value = -8500
$ws.Range("H132").Value = 5135.2856
$ws.Range("I132").Value = 2339.6956
$ws.Range("K132").Value = 7019.0868
$ws.Range("M132").Value = -4489.0868
$ws.Range("H136").Value = 20053624
$ws.Range("I136").Value = 25012944
$ws.Range("K136").Value = 75038832
$ws.Range("M136").Value = -75036282

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H37").Value = 2042
$ws.Range("I37").Value = 2042
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 2042
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -1905
$ws.Range("N37").ClearContents()
$ws.Range("H86").Value = 2725.7334
$ws.Range("I86").Value = 2647.4285
$ws.Range("K86").Value = 2647.4285
$ws.Range("M86").Value = -1524.4285
$ws.Range("H89").Value = 2725.7334
$ws.Range("I89").Value = 2647.4285
$ws.Range("K89").Value = 13237.1425
$ws.Range("M89").Value = -7621.1425
$ws.Range("H94").Value = 1748
$ws.Range("I94").Value = 2134.875
$ws.Range("J94").Value = 974.25
$ws.Range("K94").Value = 2134.875
$ws.Range("L94").Value = 974.25
$ws.Range("M94").Value = -1683.875
$ws.Range("N94").Value = -1876.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 225.625
$ws.Range("I7").Value = 242
$ws.Range("J7").Value = 198.33333
$ws.Range("K7").Value = 242
$ws.Range("L7").Value = 198.33333
$ws.Range("M7").Value = -129
$ws.Range("N7").Value = -424.33333
$ws.Range("H132").Value = 4052.5
$ws.Range("I132").Value = 2545
$ws.Range("K132").Value = 7635
$ws.Range("M132").Value = -5105

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H80").Value = 5276.6665
$ws.Range("I80").Value = 4497.5
$ws.Range("J80").Value = 5499.2856
$ws.Range("K80").Value = 13492.5
$ws.Range("L80").Value = 16497.8568
$ws.Range("M80").Value = -12556.5
$ws.Range("N80").Value = -18369.8568
$ws.Range("H83").Value = 5276.6665
$ws.Range("I83").Value = 4497.5
$ws.Range("J83").Value = 5499.2856
$ws.Range("K83").Value = 40477.5
$ws.Range("L83").Value = 49493.5704
$ws.Range("M83").Value = -35797.5
$ws.Range("N83").Value = -58853.5704
$ws.Range("H132").Value = 1638.8462
$ws.Range("J132").Value = 1637.8182
$ws.Range("L132").Value = 14740.3638
$ws.Range("N132").Value = -19800.3638
$ws.Range("H141").Value = 9341.333000000001
$ws.Range("I141").Value = 7350
$ws.Range("J141").Value = 11332.667
$ws.Range("K141").Value = 22050
$ws.Range("L141").Value = 33998.001
$ws.Range("M141").Value = -16870
$ws.Range("N141").Value = -44358.001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H116").Value = 143000
$ws.Range("J116").Value = 143000
$ws.Range("L116").Value = 143000
$ws.Range("N116").Value = -152178
$ws.Range("H122").Value = 1272.7
$ws.Range("I122").Value = 989.8333
$ws.Range("J122").Value = 1697
$ws.Range("K122").Value = 2969.4999
$ws.Range("L122").Value = 5091
$ws.Range("M122").Value = -519.4998999999998
$ws.Range("N122").Value = -9991
$ws.Range("H126").Value = 3755.9412
$ws.Range("I126").Value = 3822.6365
$ws.Range("K126").Value = 11467.9095
$ws.Range("M126").Value = -8997.9095
$ws.Range("H132").Value = 28573568
$ws.Range("I132").Value = 38463940
$ws.Range("J132").Value = 1387.7778
$ws.Range("K132").Value = 115391820
$ws.Range("L132").Value = 4163.3334
$ws.Range("M132").Value = -115389290
$ws.Range("N132").Value = -9223.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4254.0415
$ws.Range("I40").Value = 3852.2104
$ws.Range("K40").Value = 3852.2104
$ws.Range("M40").Value = -3716.2104
$ws.Range("H46").Value = 6952
$ws.Range("I46").Value = 2260.7144
$ws.Range("J46").Value = 13519.8
$ws.Range("K46").Value = 2260.7144
$ws.Range("L46").Value = 13519.8
$ws.Range("M46").Value = -2072.7144
$ws.Range("N46").Value = -13895.8
$ws.Range("H55").Value = 90909730
$ws.Range("I55").Value = 111111660
$ws.Range("K55").Value = 111111660
$ws.Range("M55").Value = -111111487
$ws.Range("H100").Value = 3866.5557
$ws.Range("I100").Value = 3600.6
$ws.Range("K100").Value = 3600.6
$ws.Range("M100").Value = -3059.6
$ws.Range("H122").Value = 6567.933
$ws.Range("I122").Value = 6293.3335
$ws.Range("K122").Value = 18880.0005
$ws.Range("M122").Value = -16430.0005
$ws.Range("H132").Value = 61354.676
$ws.Range("I132").Value = 36967.414
$ws.Range("K132").Value = 110902.242
$ws.Range("M132").Value = -108372.242

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 12498.5
$ws.Range("J45").Value = 12665.333
$ws.Range("L45").Value = 12665.333
$ws.Range("N45").Value = -13647.333
$ws.Range("H132").Value = 1940.1666
$ws.Range("I132").Value = 1710.0588
$ws.Range("J132").Value = 2499
$ws.Range("K132").Value = 5130.1764
$ws.Range("L132").Value = 7497
$ws.Range("M132").Value = -2600.1764
$ws.Range("N132").Value = -12557
